# Applies the "Creados niveles definitivos (faltan ayudas)" edit:
#  - rename the single sheet "nivel 0" -> "niveles"
#  - add a blank helper cell at P58 (" ")
#  - clear a batch of stray formatting-only cells that carried a redundant
#    "empty" style (no value) so they disappear from the sheet entirely
#  - update the window view (zoom / scroll position / selection)
#  - nudge a couple of column widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet
$ws.Name = "niveles"

# 2) New helper placeholder cell
$ws.Range("P58").Value = " "

# 3) Remove the stray empty-but-styled cells (they become fully blank/default)
$deadCells = @(
    "O4","O8",
    "I10","M10","Q10","U10",
    "G12","W12",
    "I14","U14",
    "G16","W16",
    "G20","O20","W20",
    "G24","O24","W24",
    "I26","M26","Q26","U26",
    "O28",
    "O32",
    "O35",
    "O63",
    "O66",
    "O94",
    "O97",
    "O125","O128",
    "O156","O159",
    "O187",
    "O218"
)
foreach ($addr in $deadCells) {
    $ws.Range($addr).Clear()
}

# 4) Column width tweaks (column AD a touch narrower, AE:AK widened)
$ws.Columns("AD").ColumnWidth = 1.9
$ws.Columns("AE:AK").ColumnWidth = 4.73

# 5) Window / view state: zoom, scroll position, active cell
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P59").Select() | Out-Null
